$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the regex-match repr values for the two existing rows (row 1 & 2)
$ws.Range("C1").Value = "<re.Match object; span=(8, 22), match='atest@test.com'>"
$ws.Range("C2").Value = "<re.Match object; span=(9, 23), match='atest@test.com'>"

# Append two more rows duplicating the resume rows, now with a plain email value.
# Force the "B" column to stay plain text (like the existing "29%" cells) rather
# than being auto-converted to a percentage number.
$ws.Range("A3").Value = "resume_test.docx"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "29%"
$ws.Range("C3").Value = "atest@test.com"

$ws.Range("A4").Value = "resume_test.pdf"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "29%"
$ws.Range("C4").Value = "atest@test.com"
